{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"287\u00d76=1722\", \"364\u00d78=2912\"],\n  [\"523\u00d78=4184\", \"269\u00d75=1345\"],\n  [\"318\u00d78=2544\", \"351\u00d72=702\"],\n  [\"947\u00d79=8523\", \"459\u00d73=1377\"],\n  [\"200\u00d73=600\", \"421\u00d78=3368\"],\n  [\"693\u00d76=4158\", \"326\u00d72=652\"],\n  [\"113\u00d79=1017\", \"711\u00d74=2844\"],\n  [\"670\u00d78=5360\", \"795\u00d73=2385\"],\n  [\"607\u00d76=3642\", \"837\u00d74=3348\"],\n  [\"686\u00d75=3430\", \"718\u00d77=5026\"],\n  [\"852\u00d79=7668\", \"627\u00d74=2508\"],\n  [\"559\u00d74=2236\", \"885\u00d73=2655\"],\n  [\"730\u00d79=6570\", \"676\u00d78=5408\"],\n  [\"554\u00d75=2770\", \"889\u00d78=7112\"],\n  [\"310\u00d79=2790\", \"877\u00d74=3508\"],\n  [\"133\u00d78=1064\", \"595\u00d75=2975\"],\n  [\"450\u00d77=3150\", \"511\u00d74=2044\"],\n  [\"129\u00d75=645\", \"800\u00d77=5600\"],\n  [\"509\u00d74=2036\", \"829\u00d73=2487\"],\n  [\"736\u00d74=2944\", \"639\u00d77=4473\"],\n  [\"412\u00d73=1236\", \"958\u00d72=1916\"],\n  [\"425\u00d76=2550\", \"853\u00d79=7677\"],\n  [\"916\u00d77=6412\", \"658\u00d72=1316\"],\n  [\"638\u00d75=3190\", \"329\u00d74=1316\"],\n  [\"102\u00d73=306\", \"717\u00d78=5736\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n  for (const item of searchResults.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"287\u00d76=1722\";  New = \"364\u00d78=2912\" },\n    @{ Old = \"523\u00d78=4184\";  New = \"269\u00d75=1345\" },\n    @{ Old = \"318\u00d78=2544\";  New = \"351\u00d72=702\" },\n    @{ Old = \"947\u00d79=8523\";  New = \"459\u00d73=1377\" },\n    @{ Old = \"200\u00d73=600\";   New = \"421\u00d78=3368\" },\n    @{ Old = \"693\u00d76=4158\";  New = \"326\u00d72=652\" },\n    @{ Old = \"113\u00d79=1017\";  New = \"711\u00d74=2844\" },\n    @{ Old = \"670\u00d78=5360\";  New = \"795\u00d73=2385\" },\n    @{ Old = \"607\u00d76=3642\";  New = \"837\u00d74=3348\" },\n    @{ Old = \"686\u00d75=3430\";  New = \"718\u00d77=5026\" },\n    @{ Old = \"852\u00d79=7668\";  New = \"627\u00d74=2508\" },\n    @{ Old = \"559\u00d74=2236\";  New = \"885\u00d73=2655\" },\n    @{ Old = \"730\u00d79=6570\";  New = \"676\u00d78=5408\" },\n    @{ Old = \"554\u00d75=2770\";  New = \"889\u00d78=7112\" },\n    @{ Old = \"310\u00d79=2790\";  New = \"877\u00d74=3508\" },\n    @{ Old = \"133\u00d78=1064\";  New = \"595\u00d75=2975\" },\n    @{ Old = \"450\u00d77=3150\";  New = \"511\u00d74=2044\" },\n    @{ Old = \"129\u00d75=645\";   New = \"800\u00d77=5600\" },\n    @{ Old = \"509\u00d74=2036\";  New = \"829\u00d73=2487\" },\n    @{ Old = \"736\u00d74=2944\";  New = \"639\u00d77=4473\" },\n    @{ Old = \"412\u00d73=1236\";  New = \"958\u00d72=1916\" },\n    @{ Old = \"425\u00d76=2550\";  New = \"853\u00d79=7677\" },\n    @{ Old = \"916\u00d77=6412\";  New = \"658\u00d72=1316\" },\n    @{ Old = \"638\u00d75=3190\";  New = \"329\u00d74=1316\" },\n    @{ Old = \"102\u00d73=306\";   New = \"717\u00d78=5736\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$r.Old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]2, [ref]$false, [ref]$r.New, [ref]2) | Out-Null\n}\n"}
